# Populate the distortion-coefficient rows (k00..p02, rows 6-19) for the
# cameras that previously held placeholder zeros: cam1 (col B), cam2 (col E),
# cam7 (col T), cam8 (col W). Values come from the multi-camera calibration
# processing pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cam1 (B), cam2 (E), cam7 (T), cam8 (W) results, keyed by row.
$camValues = @{
    6  = @{ B = 377.2606276170851;        E = 412.44273113994024;      T = 395.9966462212447;        W = 428.4399024684699 }
    7  = @{ B = 1.2400965644292186;       E = 1.2232746588946932;      T = 1.2152556472093268;       W = 1.2306754303475338 }
    8  = @{ B = 331.5147832330816;        E = 219.59774148411486;      T = 273.7783356609367;         W = 187.4972261282927 }
    9  = @{ B = -0.14922891541558994;     E = 0.1496739143036035;      T = -0.01619303829569439;      W = -0.07592451931818057 }
    10 = @{ B = 0.9248203156030905;       E = 0.9746196742343951;      T = 0.9907091224875562;        W = 0.9506050197585804 }
    11 = @{ B = 0.0003462762581424422;    E = -0.000358118030384714;   T = 0.00004044785203075463;    W = 0.00016911865676257648 }
    12 = @{ B = -0.00002269046045838749;  E = -0.000037350559251486395; T = -0.00005222330012631503;  W = -0.00004216212493776458 }
    13 = @{ B = -0.0000002005066477423915; E = 0.00000021238130228933882; T = -0.00000002513482382625756; W = -0.00000009285409473796876 }
    14 = @{ B = 1.2476296761807386;       E = 1.0417452308799944;      T = 1.0907434181016935;        W = 1.0612544103922412 }
    15 = @{ B = -0.0005834960968566276;   E = -0.00010207280455528885; T = -0.0002266703058971302;    W = -0.0001490832514930261 }
    16 = @{ B = 0.012837114162999583;     E = 0.041510342523125375;    T = 0.1445421001084877;        W = -0.049401322707888405 }
    17 = @{ B = 0.00000034307080697785266; E = 0.00000006212178580816405; T = 0.00000014127907221596413; W = 0.00000009042779751067861 }
    18 = @{ B = -0.000008537884274194477; E = -0.00003024484853351076; T = -0.00011561615940378931;   W = 0.000034330498960650356 }
    19 = @{ B = -0.4736463588443982;      E = -0.4372205776359586;     T = -0.4232362664922832;       W = -0.5167679930787135 }
}

foreach ($row in $camValues.Keys) {
    $cols = $camValues[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
